# "copied mainboard KiCad files" — pin-mapping refresh on the V1.1 sheet.
# The KiCad re-export swapped a couple of signal names and replaced the
# ENC1A/ENC1B/ENC2A/ENC2B encoder-pin labels with the new
# ENC1_DIR/ENC1_SPEED/ENC2_DIR/ENC2_SPEED naming scheme.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Servo outputs swapped (PA00/PA01)
$ws.Range("E4").Value = "SERVO2"
$ws.Range("E5").Value = "SERVO1"

# Encoder 1 pins (PA14/PA15) renamed
$ws.Range("E26").Value = "ENC1_SPEED"
$ws.Range("E27").Value = "ENC1_DIR"

# Encoder 2 pins (PA20/PA21) renamed
$ws.Range("E32").Value = "ENC2_DIR"
$ws.Range("E33").Value = "ENC2_SPEED"

# Motor 2 outputs swapped (PA22/PA23)
$ws.Range("E34").Value = "MOTOR2B"
$ws.Range("E35").Value = "MOTOR2A"

# Leave the selection where the author's session ended up.
$ws.Range("E36").Select()
